# Weekly update: insert a new "Espinaca" price record for Vega Central
# Mapocho de Santiago (Provincia de Chacabuco) at row 464, pushing the
# existing historical rows (old 464-497) down by one row (new 465-498).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 464; this shifts rows 464:497 down
# to 465:498 and extends the sheet dimension accordingly.
$ws.Rows.Item(464).Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A464").Value2 = 9
$ws.Range("B464").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C464").Value2 = "Metropolitana"
$ws.Range("D464").Value2 = 44826
$ws.Range("E464").Value2 = 13
$ws.Range("F464").Value2 = 100112012
$ws.Range("G464").Value2 = "Espinaca"
$ws.Range("H464").Value2 = "Sin especificar"
$ws.Range("I464").Value2 = "Primera"
$ws.Range("J464").Value2 = 160
$ws.Range("K464").Value2 = 6000
$ws.Range("L464").Value2 = 7000
$ws.Range("M464").Value2 = 6500
$ws.Range("N464").Value2 = "`$/cuna 10 kilos"
$ws.Range("O464").Value2 = "Provincia de Chacabuco"
$ws.Range("P464").Value2 = 650
$ws.Range("Q464").Value2 = 10
$ws.Range("R464").Value2 = "Hortaliza"

# Ensure the date cell keeps the same date/time number format used by
# the other rows in column D.
$ws.Range("D464").NumberFormat = $ws.Range("D465").NumberFormat
